$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Per-row Price (D) / Volume(1h) (E) updates ---
# Note: some Price values look like plain decimal numbers (e.g. "513.34").
# The source data stores Price as TEXT (not a number), so for any new value that
# Excel would auto-convert to a number on assignment, we briefly force the cell to
# Text format, assign the literal string, then restore the cell style to "Normal"
# so no stray number-format styling is left behind (matches original, unstyled cells).
$ws.Range("D2").Value = "57.051.31"
$ws.Range("E2").Value = "  +2.59%  "

$ws.Range("D3").Value = "3.011.54"
$ws.Range("E3").Value = "  +1.59%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.437"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.50%  "

$ws.Range("E10").Value = "  +7.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.359"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("E12").Value = "  +1.84%  "

$ws.Range("D13").Value = "3.524.03"
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.09%  "

$ws.Range("E15").Value = "  +11.65%  "

$ws.Range("D16").Value = "57.034.16"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").Value = "3.008.16"
$ws.Range("E17").Value = "  +1.71%  "

$ws.Range("E18").Value = "  +5.15%  "

$ws.Range("E19").Value = "  +3.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.15%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.489"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.172"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.66%  "

$ws.Range("E26").Value = "  -2.60%  "

$ws.Range("D27").Value = "0.0₃0920"
$ws.Range("E27").Value = "  +8.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.37%  "

$ws.Range("E30").Value = "  +5.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.79%  "

$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0682"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.66%  "

$ws.Range("E38").Value = "  +3.31%  "

$ws.Range("D39").Value = "3.042.64"
$ws.Range("E39").Value = "  +1.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.24%  "

$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").Value = "2.301.24"
$ws.Range("E42").Value = "  +7.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.649"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.67%  "

$ws.Range("E46").Value = "  +3.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0875"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.78%  "

# Row 44 & 45: Filecoin and ONDO swap rank positions (with updated Price/Volume)
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.40%  "

# Row 48 & 49: Cosmos and VeChain swap rank positions (with updated Price/Volume)
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0240"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.23%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.04%  "
